$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  2 = @{ D='27.934.71'; E='  +1.59%  ' }
  3 = @{ D='1.636.45'; E='  +0.33%  ' }
  4 = @{ D='0.998'; E='  -0.27%  ' }
  5 = @{ D='212.48'; E='  +0.48%  ' }
  6 = @{ D='0.522'; E='  +0.44%  ' }
  7 = @{ D='0.998'; E='  -0.26%  ' }
  8 = @{ D='23.40'; E='  +1.48%  ' }
  9 = @{ D='0.258'; E='  -1.72%  ' }
  10 = @{ D='0.0612'; E='  +0.36%  ' }
  11 = @{ D='0.0882'; E='  +2.50%  ' }
  12 = @{ D='1.863.29'; E='  +0.13%  ' }
  13 = @{ D='1.635.10'; E='  +0.21%  ' }
  14 = @{ D='4.07'; E='  +0.94%  ' }
  15 = @{ D='0.571'; E='  +2.45%  ' }
  16 = @{ D='65.37'; E='  +0.65%  ' }
  17 = @{ D='27.885.41'; E='  +1.05%  ' }
  18 = @{ D='232.44'; E='  +1.23%  ' }
  19 = @{ D='0.0₃0720'; E='  +0.42%  ' }
  20 = @{ D='7.58'; E='  +0.55%  ' }
  21 = @{ D='0.996'; E='  -0.31%  ' }
  22 = @{ D='10.50'; E='  -2.39%  ' }
  23 = @{ D='4.36'; E='  -0.09%  ' }
  24 = @{ D='2.10'; E='  -1.35%  ' }
  25 = @{ D='152.13'; E='  +1.85%  ' }
  26 = @{ D='6.87'; E='  +0.19%  ' }
  27 = @{ D='15.65'; E='  +0.46%  ' }
  28 = @{ E='  +0.02%  ' }
  29 = @{ D='0.997'; E='  -0.23%  ' }
  30 = @{ E='  +0.10%  ' }
  31 = @{ D='0.0483'; E='  +0.45%  ' }
  32 = @{ E='  +2.34%  ' }
  33 = @{ D='3.11'; E='  +0.69%  ' }
  34 = @{ D='1.405.22'; E='  -4.03%  ' }
  35 = @{ E='  +2.33%  ' }
  36 = @{ D='2.35'; E='  +1.30%  ' }
  37 = @{ E='  +1.17%  ' }
  38 = @{ D='0.879'; E='  +0.52%  ' }
  39 = @{ D='0.559'; E='  +0.15%  ' }
  40 = @{ D='0.917'; E='  +0.58%  ' }
  41 = @{ E='  +1.10%  ' }
  42 = @{ D='0.997'; E='  -0.25%  ' }
  43 = @{ D='67.32'; E='  -0.81%  ' }
  44 = @{ B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.29'; E='  +0.48%  ' }
  45 = @{ B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='1.85'; E='  +6.23%  ' }
  46 = @{ B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='5.49'; E='  +2.42%  ' }
  47 = @{ D='1.773.66'; E='  +0.40%  ' }
  48 = @{ D='87.59'; E='  +0.34%  ' }
  49 = @{ D='0.0999'; E='  +0.45%  ' }
  50 = @{ D='0.0506'; E='  +0.40%  ' }
  51 = @{ D='7.61'; E='  -1.26%  ' }
}

foreach ($rowNum in $updates.Keys) {
  $rowData = $updates[$rowNum]
  foreach ($col in $rowData.Keys) {
    $cellRef = "$col$rowNum"
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.Value = "'" + $rowData[$col]
    $cell.Style = $origStyle
  }
}

Write-Host "Update complete"